$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Copy formatting for the new rows (57-60) from the last existing data row (56) ---
$ws.Range("A56:D56").Copy() | Out-Null
$ws.Range("A57:D57").PasteSpecial(-4122) | Out-Null
$ws.Range("A56:D56").Copy() | Out-Null
$ws.Range("A58:D58").PasteSpecial(-4122) | Out-Null
$ws.Range("A56:D56").Copy() | Out-Null
$ws.Range("A59:D59").PasteSpecial(-4122) | Out-Null
$ws.Range("A56:D56").Copy() | Out-Null
$ws.Range("A60:D60").PasteSpecial(-4122) | Out-Null

$ws.Rows("57").RowHeight = 15.75
$ws.Rows("58").RowHeight = 15.75
$ws.Rows("59").RowHeight = 15.75
$ws.Rows("60").RowHeight = 15.75

# --- New logbook entries ---
$ws.Range("A57").Value2 = 22
$ws.Range("B57").Value2 = 1.5

$ws.Range("A58").Value2 = 23
$ws.Range("B58").Value2 = 4

$ws.Range("A59").Value2 = 24
$ws.Range("B59").Value2 = 6.5

$ws.Range("A60").Value2 = 25
$ws.Range("B60").Value2 = 6

# Descriptions - written in an order that matches the original shared-string table layout
$ws.Range("D59").Value2 = "Layout update: custom checkboxes / buttons en andere dingen"
$ws.Range("D57").Value2 = "Code bestuderen en bugfixen"
$ws.Range("D60").Value2 = "Skype + verder werken layout + rivieren in nieuwe kaart verwerken"
$ws.Range("D58").Value2 = "Verder werken aan app"

# --- Cumulative-hours formulas ---
# Existing block (rows 25-38) re-entered as one range so it collapses into a shared formula group
$ws.Range("C25:C38").Formula = "=SUM(B25,C24)"
# New rows cumulative total, continuing the C44:C56 pattern
$ws.Range("C57:C60").Formula = "=(SUM(C56,B57))"

# --- Grand total now also includes the new last row ---
$ws.Range("H1").Formula = "=(SUM(C20,C38,C60))"

# --- View state: scrolled/selected like after entering the new rows ---
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("B60").Select() | Out-Null
